{"js": "// Enhance Siege Analytics descriptions with voter file discovery metrics.\n//\n// 1. Insert three new bullet paragraphs immediately before the\n//    \"Developed and deployed custom analytical tools...\" bullet.\n// 2. Remove the \"Created fraud detection systems for campaign finance\n//    data analysis across multi-terabyte datasets\" bullet (it followed the\n//    \"...170% more viable targets\" bullet).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText =\n  \"\u2022 Developed and deployed custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\";\nconst removeText =\n  \"\u2022 Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets\";\n\nconst newBullets = [\n  \"\u2022 Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data\",\n  \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"\u2022 Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts\",\n];\n\nlet anchorParagraph = null;\nlet removeParagraph = null;\n\nfor (const p of paragraphs.items) {\n  const text = p.text.trim();\n  if (anchorParagraph === null && text === anchorText) {\n    anchorParagraph = p;\n  }\n  if (removeParagraph === null && text === removeText) {\n    removeParagraph = p;\n  }\n}\n\nif (anchorParagraph) {\n  // insertParagraph(text, \"Before\") inserts each new paragraph directly\n  // above the anchor, so insert in order to preserve the intended sequence.\n  for (const bulletText of newBullets) {\n    anchorParagraph.insertParagraph(bulletText, \"Before\");\n  }\n}\n\nif (removeParagraph) {\n  removeParagraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Enhance Siege Analytics descriptions with voter file discovery metrics.\n#\n# 1. Insert three new bullet paragraphs immediately before the\n#    \"Developed and deployed custom analytical tools...\" bullet.\n# 2. Remove the \"Created fraud detection systems for campaign finance\n#    data analysis across multi-terabyte datasets\" bullet (it followed the\n#    \"...170% more viable targets\" bullet).\n\n$d = $word.ActiveDocument\n$cr = [char]13\n\n$anchorText = \"\u2022 Developed and deployed custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\"\n$removeText = \"\u2022 Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets\"\n\n$newBullets = \"\u2022 Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data\" + $cr + `\n    \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\" + $cr + `\n    \"\u2022 Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts\" + $cr\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $anchorText) {\n        $p.Range.InsertBefore($newBullets)\n        break\n    }\n}\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $removeText) {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
